# Actualizo con datos fiscales y de comex para diciembre22
$wb = $excel.ActiveWorkbook

# --- "BC por zonas": trade balance by zone (Mercosur, Chile, ALADI, ...) ---
$wsBC = $wb.Worksheets.Item("BC por zonas")

$bcData = @{
    2 = @(15793, 19264)
    3 = @(4938, 778)
    4 = @(5313, 3377)
    5 = @(1480, 35)
    6 = @(8653, 12557)
    7 = @(10846, 11118)
    8 = @(685, 652)
    9 = @(1092, 586)
    10 = @(581, 432)
    11 = @(7894, 4225)
    12 = @(8022, 17516)
    13 = @(2020, 729)
    14 = @(795, 1201)
    15 = @(4555, 1849)
    16 = @(4655, 2534)
    17 = @(3964, 1179)
    18 = @(463, 195)
    19 = @(864, 515)
    20 = @(5832, 2779)
}

foreach ($row in $bcData.Keys) {
    $vals = $bcData[$row]
    $wsBC.Range("B$row").Value = $vals[0]
    $wsBC.Range("C$row").Value = $vals[1]
    $wsBC.Range("D$row").Formula = "=B$row-C$row"
}

# Column D (saldo) now shows thousands separators (previously General format)
$wsBC.Range("D2:D20").NumberFormat = "#,##0"

# Drop the border that used to separate the last two rows from the rest of the table
$wsBC.Range("B19:C20").Borders.LineStyle = -4142

# --- "Expo-ICA": exports by country/region (Comex, Dec-22) ---
$wsExpo = $wb.Worksheets.Item("Expo-ICA")
$expoData = @{
    2 = 88446
    3 = 23868
    4 = 28
    5 = 1539
    6 = 248
    7 = 635
    8 = 525
    9 = 15575
    10 = 4347
    11 = 317
    12 = 28
    13 = 177
    14 = 300
    15 = 148
    16 = 33119
    17 = 4158
    18 = 265
    19 = 1422
    20 = 58
    21 = 122
    22 = 177
    23 = 1020
    24 = 9170
    25 = 297
    26 = 814
    27 = 948
    28 = 13249
    29 = 319
    30 = 463
    31 = 127
    32 = 510
    33 = 23061
    34 = 6119
    35 = 1083
    36 = 225
    37 = 30
    38 = 464
    39 = 143
    40 = 9
    41 = 177
    42 = 2680
    43 = 2043
    44 = 1594
    45 = 7950
    46 = 247
    47 = 296
    48 = 8398
    49 = 3867
    50 = 2573
    51 = 90
    52 = 1369
    53 = 498
}
foreach ($row in $expoData.Keys) {
    $wsExpo.Range("B$row").Value = $expoData[$row]
}

# --- "Impo-ICA": imports by country/region (Comex, Dec-22) ---
$wsImpo = $wb.Worksheets.Item("Impo-ICA")
$impoData = @{
    2 = 81523
    3 = 12454
    4 = 9352
    5 = 1340
    6 = 1318
    7 = 444
    8 = 30009
    9 = 2533
    10 = 1461
    11 = 12381
    12 = 3842
    13 = 1199
    14 = 1388
    15 = 671
    16 = 4329
    17 = 2204
    18 = 12868
    19 = 12500
    20 = 368
    21 = 15037
    22 = 9705
    23 = 4231
    24 = 1100
    25 = 8567
    26 = 280
    27 = 629
    28 = 790
    29 = 2538
    30 = 512
    31 = 423
    32 = 697
    33 = 628
    34 = 578
    35 = 316
    36 = 569
    37 = 605
    38 = 1996
    39 = 1996
    40 = 592
}
foreach ($row in $impoData.Keys) {
    $wsImpo.Range("B$row").Value = $impoData[$row]
}

# --- Update selections left behind by the editing session, and move the active tab ---
$wsBC.Range("I15").Select()
$wsImpo.Range("B2:B40").Select()

